$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '29.094.36'
$ws.Cells.Item(2, 5).Value = '  +0.21%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.833.80'
$ws.Cells.Item(3, 5).Value = '  +0.05%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.41%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '243.69'
$ws.Cells.Item(5, 5).Value = '  +0.85%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.6287'
$ws.Cells.Item(6, 5).Value = '  +0.10%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '1.003'
$ws.Cells.Item(7, 5).Value = '  +0.30%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.07462'
$ws.Cells.Item(8, 5).Value = '  -2.00%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.2928'
$ws.Cells.Item(9, 5).Value = '  +0.31%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.46%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07723'
$ws.Cells.Item(11, 5).Value = '  -0.22%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.837.91'
$ws.Cells.Item(12, 5).Value = '  +0.26%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.996'
$ws.Cells.Item(13, 5).Value = '  +0.72%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.6679'
$ws.Cells.Item(14, 5).Value = '  +0.45%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '82.91'
$ws.Cells.Item(15, 5).Value = '  +0.06%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.000009347'
$ws.Cells.Item(16, 5).Value = '  -4.30%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '6.068'
$ws.Cells.Item(17, 5).Value = '  +1.09%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '29.110.01'
$ws.Cells.Item(18, 5).Value = '  +0.31%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +2.23%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '223.56'
$ws.Cells.Item(20, 5).Value = '  -1.09%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '1.004'
$ws.Cells.Item(21, 5).Value = '  +0.47%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '7.127'
$ws.Cells.Item(22, 5).Value = '  -1.02%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.35%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '160.03'
$ws.Cells.Item(24, 5).Value = '  +1.12%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.1402'
$ws.Cells.Item(25, 5).Value = '  +2.25%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.509'
$ws.Cells.Item(26, 5).Value = '  +0.98%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '17.91'
$ws.Cells.Item(27, 5).Value = '  +0.15%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.498'
$ws.Cells.Item(28, 5).Value = '  +0.47%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '4.146'
$ws.Cells.Item(29, 5).Value = '  +1.97%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.09%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.05477'
$ws.Cells.Item(31, 5).Value = '  +5.39%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.43%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.7516'
$ws.Cells.Item(33, 5).Value = '  +1.66%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.852'
$ws.Cells.Item(34, 5).Value = '  +0.28%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.135'
$ws.Cells.Item(35, 5).Value = '  -0.81%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.607'
$ws.Cells.Item(36, 5).Value = '  -3.30%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.228.86'
$ws.Cells.Item(37, 5).Value = '  -3.14%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.753'
$ws.Cells.Item(38, 5).Value = '  -0.23%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01785'
$ws.Cells.Item(39, 5).Value = '  -0.17%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '6.643'
$ws.Cells.Item(40, 5).Value = '  +6.23%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.8943'
$ws.Cells.Item(41, 5).Value = '  -0.07%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.003'
$ws.Cells.Item(42, 5).Value = '  +0.33%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '101.69'
$ws.Cells.Item(43, 5).Value = '  +0.12%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '65.60'
$ws.Cells.Item(44, 5).Value = '  +1.47%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.00000000125'
$ws.Cells.Item(45, 5).Value = '  +1.70%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'XinFinNetwork'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.07728'
$ws.Cells.Item(46, 5).Value = '  +10.53%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Mantle'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.5100'
$ws.Cells.Item(47, 5).Value = '  -0.22%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.4044'
$ws.Cells.Item(48, 5).Value = '  +1.32%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '9.002'
$ws.Cells.Item(49, 5).Value = '  +1.78%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.05806'
$ws.Cells.Item(50, 5).Value = '  +0.89%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'RenderToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.656'
$ws.Cells.Item(51, 5).Value = '  +1.60%  '
